$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Attendance")

# "Normal Curfew" -> "Leaving Camp Curfew"
$ws.Range("A2").Value = "Leaving Camp Curfew"

# Update the active selection to match the merged cell region A3:G3
$ws.Activate()
$ws.Range("A3:G3").Select()
